# 21.turnover_worker_changed.xlsx
# Commit message: "translate the formula downwards"
#
# The 粗利 (profit) column D6:D15 was blank; the SUM formula in D3 sums
# D6:D15 downward, and each row's profit (売価 - 原価, i.e. column C - column B)
# gets filled in as the worked/translated-down result. The sheet is also
# renamed from the month label "4月" to the "yyyymm" form "202004", the
# B:D columns are resized to a uniform width, and the active selection
# moves to the summary cell C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from "4月" to "202004"
$ws.Name = "202004"

# Translate/fill the (sale price - cost) result down column D, row by row
$ws.Range("D6").Value = 126000
$ws.Range("D7").Value = 85000
$ws.Range("D8").Value = 91000
$ws.Range("D9").Value = 132000
$ws.Range("D10").Value = 137600
$ws.Range("D11").Value = 124000
$ws.Range("D12").Value = 133600
$ws.Range("D13").Value = 136800
$ws.Range("D14").Value = 133600
$ws.Range("D15").Value = 136800

# Resize columns B:D to a uniform width (stored width == 11)
$ws.Columns("B:D").ColumnWidth = 10.285714285714286

# Move the selection to the summary formula cell
$ws.Range("C3").Select()
